$d = $word.ActiveDocument

# 1. Update the letter date: September 19, 2025 -> September 21, 2025
$null = $d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line "3037 Lamory Pl, Santa Clara CA 95051"
#    (the one in the letter body, not the one inside the info table) into
#    two paragraphs: "3037 Lamory Pl" and a new "Santa Clara, CA 95051".
$count = $d.Paragraphs.Count
$addrIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13)
    $inTable = $para.Range.Information(12)
    if ((-not $inTable) -and ($txt -eq "3037 Lamory Pl, Santa Clara CA 95051")) {
        $addrIndex = $i
        break
    }
}

if ($addrIndex -gt 0) {
    $addrPara = $d.Paragraphs.Item($addrIndex)
    $addrPara.Range.Text = "3037 Lamory Pl"

    $addrPara2 = $d.Paragraphs.Item($addrIndex)
    $null = $addrPara2.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($addrIndex + 1)
    $newPara.Range.Text = "Santa Clara, CA 95051"
}

# 3. Remove the empty "No Spacing" paragraph that immediately follows the
#    "Board of Directors" signature line.
$count2 = $d.Paragraphs.Count
for ($i = 1; $i -le $count2; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.TrimEnd([char]13)
    if ($txt -match "Board of Directors$") {
        $next = $d.Paragraphs.Item($i + 1)
        $ntxt = $next.Range.Text.TrimEnd([char]13)
        if ($ntxt -eq "" -and $next.Range.ParagraphStyle.NameLocal -eq "No Spacing") {
            $null = $next.Range.Delete()
        }
        break
    }
}

Write-Output "edits applied"
